$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '63.288.94' },
    @{ Cell = 'E2'; Value = '  +1.06%  ' },
    @{ Cell = 'D3'; Value = '2.651.38' },
    @{ Cell = 'E3'; Value = '  +3.41%  ' },
    @{ Cell = 'D4'; Value = '0.999' },
    @{ Cell = 'E4'; Value = '  -0.06%  ' },
    @{ Cell = 'D5'; Value = '594.58' },
    @{ Cell = 'E5'; Value = '  +2.80%  ' },
    @{ Cell = 'D6'; Value = '143.66' },
    @{ Cell = 'E6'; Value = '  +0.26%  ' },
    @{ Cell = 'D7'; Value = '0.999' },
    @{ Cell = 'E7'; Value = '  -0.04%  ' },
    @{ Cell = 'E8'; Value = '  +0.28%  ' },
    @{ Cell = 'D9'; Value = '2.647.48' },
    @{ Cell = 'E9'; Value = '  +3.27%  ' },
    @{ Cell = 'E10'; Value = '  +0.66%  ' },
    @{ Cell = 'D12'; Value = '0.153' },
    @{ Cell = 'E12'; Value = '  +0.87%  ' },
    @{ Cell = 'E13'; Value = '  +1.89%  ' },
    @{ Cell = 'D14'; Value = '27.46' },
    @{ Cell = 'E14'; Value = '  +2.65%  ' },
    @{ Cell = 'D15'; Value = '3.123.36' },
    @{ Cell = 'E15'; Value = '  +3.33%  ' },
    @{ Cell = 'D16'; Value = '63.183.86' },
    @{ Cell = 'E16'; Value = '  +1.05%  ' },
    @{ Cell = 'E17'; Value = '  +0.59%  ' },
    @{ Cell = 'D18'; Value = '2.651.98' },
    @{ Cell = 'E18'; Value = '  +3.33%  ' },
    @{ Cell = 'D19'; Value = '11.46' },
    @{ Cell = 'E19'; Value = '  +3.45%  ' },
    @{ Cell = 'D20'; Value = '340.04' },
    @{ Cell = 'E20'; Value = '  +0.71%  ' },
    @{ Cell = 'D21'; Value = '4.38' },
    @{ Cell = 'E21'; Value = '  +1.48%  ' },
    @{ Cell = 'D22'; Value = '6.80' },
    @{ Cell = 'E22'; Value = '  +2.44%  ' },
    @{ Cell = 'E23'; Value = '  +0.02%  ' },
    @{ Cell = 'D24'; Value = '67.35' },
    @{ Cell = 'E24'; Value = '  +0.69%  ' },
    @{ Cell = 'D25'; Value = '1.69' },
    @{ Cell = 'E25'; Value = '  +7.26%  ' },
    @{ Cell = 'E26'; Value = '  +1.77%  ' },
    @{ Cell = 'E27'; Value = '  +1.16%  ' },
    @{ Cell = 'B28'; Value = 'InternetComputer(DFINITY)' },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' },
    @{ Cell = 'D28'; Value = '8.46' },
    @{ Cell = 'E28'; Value = '  +3.53%  ' },
    @{ Cell = 'B29'; Value = 'Binance-PegBSC-USD' },
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd' },
    @{ Cell = 'D29'; Value = '1.00' },
    @{ Cell = 'E29'; Value = '  -0.03%  ' },
    @{ Cell = 'D30'; Value = '536.84' },
    @{ Cell = 'E30'; Value = '  +17.80%  ' },
    @{ Cell = 'D31'; Value = '7.85' },
    @{ Cell = 'E31'; Value = '  -0.82%  ' },
    @{ Cell = 'D32'; Value = '1.85' },
    @{ Cell = 'E32'; Value = '  +14.45%  ' },
    @{ Cell = 'E33'; Value = '  +3.97%  ' },
    @{ Cell = 'D34'; Value = '0.0₃0809' },
    @{ Cell = 'E34'; Value = '  +2.22%  ' },
    @{ Cell = 'D35'; Value = '173.76' },
    @{ Cell = 'E35'; Value = '  -1.73%  ' },
    @{ Cell = 'D36'; Value = '5.11' },
    @{ Cell = 'E36'; Value = '  +15.39%  ' },
    @{ Cell = 'D37'; Value = '0.407' },
    @{ Cell = 'E37'; Value = '  +3.24%  ' },
    @{ Cell = 'D38'; Value = '0.998' },
    @{ Cell = 'E38'; Value = '  -0.26%  ' },
    @{ Cell = 'D39'; Value = '19.06' },
    @{ Cell = 'E39'; Value = '  +1.53%  ' },
    @{ Cell = 'D40'; Value = '1.85' },
    @{ Cell = 'E40'; Value = '  +9.93%  ' },
    @{ Cell = 'D41'; Value = '172.75' },
    @{ Cell = 'E41'; Value = '  +10.06%  ' },
    @{ Cell = 'D42'; Value = '0.998' },
    @{ Cell = 'E42'; Value = '  -0.16%  ' },
    @{ Cell = 'D43'; Value = '40.11' },
    @{ Cell = 'E43'; Value = '  -0.92%  ' },
    @{ Cell = 'D44'; Value = '3.76' },
    @{ Cell = 'E44'; Value = '  +2.55%  ' },
    @{ Cell = 'D45'; Value = '22.08' },
    @{ Cell = 'E45'; Value = '  +6.41%  ' },
    @{ Cell = 'D46'; Value = '0.0561' },
    @{ Cell = 'E46'; Value = '  +5.62%  ' },
    @{ Cell = 'D47'; Value = '0.633' },
    @{ Cell = 'E47'; Value = '  +0.93%  ' },
    @{ Cell = 'D48'; Value = '0.0241' },
    @{ Cell = 'E48'; Value = '  +3.38%  ' },
    @{ Cell = 'D49'; Value = '0.0961' },
    @{ Cell = 'E49'; Value = '  +0.69%  ' },
    @{ Cell = 'D50'; Value = '18.71' },
    @{ Cell = 'E50'; Value = '  +4.53%  ' },
    @{ Cell = 'D51'; Value = '1.72' },
    @{ Cell = 'E51'; Value = '  +3.54%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}

Write-Output ("Applied {0} cell updates" -f $updates.Count)
